$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 79, shifting rows 79:204 down to 80:205
$ws.Rows("79:79").Insert()

# Populate the new row 79 with its data values
$ws.Cells.Item(79, 1).Value = 5
$ws.Cells.Item(79, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(79, 3).Value = "Maule"
$ws.Cells.Item(79, 4).Value = 45162
$ws.Cells.Item(79, 4).NumberFormat = $ws.Cells.Item(80, 4).NumberFormat
$ws.Cells.Item(79, 5).Value = 7
$ws.Cells.Item(79, 6).Value = 100112001
$ws.Cells.Item(79, 7).Value = "Berenjena"
$ws.Cells.Item(79, 8).Value = "Sin especificar"
$ws.Cells.Item(79, 9).Value = "Primera"
$ws.Cells.Item(79, 10).Value = 200
$ws.Cells.Item(79, 11).Value = 10000
$ws.Cells.Item(79, 12).Value = 10000
$ws.Cells.Item(79, 13).Value = 10000
$ws.Cells.Item(79, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(79, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(79, 16).Value = 200
$ws.Cells.Item(79, 17).Value = 50
$ws.Cells.Item(79, 18).Value = "Hortaliza"
